$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for new rows 286-289 (continuing the pattern of row 285)
$dates = @(45842, 45843, 45844, 45845)
$values = @(116.4121952, 0.00170247, 0.008850780000000001, 0.06933635, 12792.90181321, 465.80531254, 0.24, 1.7904431, 485.38834923)

$startRow = 286

for ($i = 0; $i -lt $dates.Length; $i++) {
    $r = $startRow + $i

    # Copy style from the last existing data row (A285) into the new date cell (column A)
    # so it gets the same style index (date format) instead of a brand-new style.
    $ws.Cells.Item(285, 1).Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)
    $ws.Cells.Item($r, 1).Value2 = $dates[$i]

    for ($c = 0; $c -lt $values.Length; $c++) {
        $ws.Cells.Item($r, $c + 2).Value2 = $values[$c]
    }
}

$excel.CutCopyMode = 0
